$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:K data to E:L
$ws.Columns("D").Insert()

# Copy cell formatting (number format/font/alignment) from the now-shifted
# column E into the newly inserted column D so the new cells match the
# rest of the data table's styling.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows that never had data in column D/E to begin with (section headers /
# blank separators) should stay empty in column D too.
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()

# Populate the new column D with the latest period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 61300
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = -200
$ws.Range("D17").Value = 9400
$ws.Range("D18").Value = 51900
$ws.Range("D20").Value = -25900
$ws.Range("D21").Value = 29600
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 26000
$ws.Range("D24").Value = 6700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 19300
$ws.Range("D27").Value = 19300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 25900
$ws.Range("D33").Value = 19300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 19300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 13800
$ws.Range("D42").Value = 92400
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 25200
$ws.Range("D49").Value = 100
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 6200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1585300
$ws.Range("D57").Value = 14100
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1398200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 56200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 187200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 19300
$ws.Range("D83").Value = 3600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 20700
$ws.Range("D91").Value = -1600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 8100
$ws.Range("D96").Value = -6400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -58100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -29400

"done"
